# Add remote server shortcuts (show, upload, sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B22 ("U" row, Ctrl column): "Navigate to super" -> "Navigate to super / Upload to remote"
$r = $ws.Range("B22")
$r.Value = "Navigate to super / Upload to remote"
$c1 = $r.Characters(1, 18)
$c1.Font.Color = 8421504
$c2 = $r.Characters(19, 19)
$c2.Font.Color = 0

# --- C19 ("R" row, Ctrl+Shift column): "- / Replace in Files" -> "Replace in Files / Open Remote Hosts"
$r = $ws.Range("C19")
$r.Value = "Replace in Files / Open Remote Hosts"
$c1 = $r.Characters(1, 17)
$c1.Font.Color = 8421504
$c2 = $r.Characters(18, 19)
$c2.Font.Color = 0

# --- C26 ("Y" row, Ctrl+Alt column): "-" -> "- / Sync with remote"
$r = $ws.Range("C26")
$r.Value = "- / Sync with remote"
$c2 = $r.Characters(3, 18)
$c2.Font.Color = 0

# --- sheet view: active cell / top-left cell moved
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()

# --- column C width changed from 33 to 34.85546875
$ws.Columns("C").ColumnWidth = 34.85546875
